# added 4wk low sales check
$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary = $wb.Worksheets.Item("Summary")

# Update Seasonality Index (column L) on "Forecast Comparison" sheet
$wsForecast.Range("L2").Value = 1.09
$wsForecast.Range("L3").Value = 0.87
$wsForecast.Range("L5").Value = 0.95
$wsForecast.Range("L6").Value = 0.86
$wsForecast.Range("L7").Value = 1.02
$wsForecast.Range("L8").Value = 1.16
$wsForecast.Range("L9").Value = 0.89
$wsForecast.Range("L10").Value = 0.89
$wsForecast.Range("L11").Value = 0.96
$wsForecast.Range("L12").Value = 1.17
$wsForecast.Range("L13").Value = 0.93
$wsForecast.Range("L14").Value = 1.02
$wsForecast.Range("L15").Value = 1.16
$wsForecast.Range("L16").Value = 0.94
$wsForecast.Range("L17").Value = 1

# Update forecast totals on "Summary" sheet (4wk low sales check values
# changed -- these are stored as text, so keep them as text using a
# leading apostrophe, same as Excel's "number stored as text" behavior)
$wsSummary.Range("B9").Value = "'7"
$wsSummary.Range("B10").Value = "'3"
$wsSummary.Range("B11").Value = "'2"
$wsSummary.Range("B12").Value = "'0"
$wsSummary.Range("B14").Value = "'0"
